$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Finish row 70 (time logged, "?" minutes, and notes) ---
# C70 needs the same "time" style already used elsewhere in column C (s="4")
$ws.Cells.Item(69, 3).Copy() | Out-Null
$ws.Cells.Item(70, 3).PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Cells.Item(70, 3).Value = "9:34PM"

$ws.Cells.Item(70, 5).Value = "?"
$ws.Cells.Item(70, 7).Value = "Got newInvoiceCandS gui almost all working"

# --- Add new log entry row 72 ---
# A72 needs the same "date" style used by the other date cells in column A (s="2")
$ws.Cells.Item(66, 1).Copy() | Out-Null
$ws.Cells.Item(72, 1).PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Cells.Item(72, 1).Value = 43811

$ws.Cells.Item(72, 2).Value = "11:40AM"

# C72 needs the "time" style too (s="4")
$ws.Cells.Item(69, 3).Copy() | Out-Null
$ws.Cells.Item(72, 3).PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Cells.Item(72, 3).Value = "12:58PM"

$ws.Cells.Item(72, 4).Value = 10
$ws.Cells.Item(72, 5).Value = 68
$ws.Cells.Item(72, 6).Value = "Code"
$ws.Cells.Item(72, 7).Value = "Adding products to show up on newInvoiceCandS, and then beginning to code the mailto: "
$ws.Cells.Item(72, 8).Value = "GUI now correctly shows all products through a convoluted SQLite3 query. Much time taken to build the query into 1 call"

$excel.CutCopyMode = 0
$ws.Range("H72").Select() | Out-Null
